$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (LEDs): duration/period changed ---
$ws.Range("C3").Value = 16

# --- Row 5 (Voltage): duration changed ---
$ws.Range("B5").Value = 2000

# --- Row 6 (Idle): no longer derived from C6*1000*E6; B6/C6 cleared ---
$ws.Range("B6").ClearContents()
$ws.Range("C6").ClearContents()

# --- Comments column (G): update text on the Idle row; engine drops the
#     now-unused shared string and appends the new one, matching the diff ---
$ws.Range("G6").Value = "PWR_DWN; T2,ADC,SPI,WDT,BOD enabled."

# --- Battery table: make the current-draw formula a shared formula across
#     C11:C18 (previously each row carried its own un-shared copy) ---
$ws.Range("C11:C18").Formula = "=B11/AVG_CURRENT"

# --- Selection moved to G3 ---
$ws.Range("G3").Select()
